$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a new column CK that mirrors column CJ's formatting (header style +
# centered integer style for the data rows), then fill in the new
# "13-nov" day column.
$ws.Range("CJ1:CJ11").Copy()
$ws.Range("CK1:CK11").PasteSpecial(-4122)

$ws.Range("CK1").Value = "13-nov"
$ws.Range("CK2").Value = 10
$ws.Range("CK3").Value = 10
$ws.Range("CK4").Value = 8
$ws.Range("CK5").Value = 9
$ws.Range("CK6").Value = 9
$ws.Range("CK7").Value = 5
$ws.Range("CK8").Value = 15
$ws.Range("CK9").Value = 12
$ws.Range("CK10").Value = 14
$ws.Range("CK11").Value = 0

$ws.Range("CK11").Select()
